$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $true, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# Title (appears twice: Heading1 at top, and bold paragraph near the bottom)
Replace-Text "Play Miss Red Slot Game for Free - IGT Video Slots Review" "Play Miss Red Slot for Free"

# "What we like" bullet list
Replace-Text "1024 pay lines through MultiwayXtra structure" "Immersive fairy tale theme"
Replace-Text "Expanding wolf and Miss Red symbols" "Top-quality graphics and sound"
Replace-Text "Free spins with potential for reactivation" "MultiwayXtra winning structure"
Replace-Text "Top-quality graphics and sound design" "Opportunity to reactivate free spins round"

# "What we don't like" bullet list
Replace-Text "Miss Red's sexualized appearance may not be for everyone" "Sensualized depiction of Little Red Riding Hood"
Replace-Text "Medium-level variance may not appeal to high-stakes players" "Limited betting options"

# Meta description (italic paragraph)
Replace-Text "Read our review of Miss Red online slot game by IGT. Play for free and find out about the game's MultiwayXtra structure, expanding wilds, and bonus features." "Experience the immersive fairy tale-themed game and enjoy the chance to play for free with Miss Red."
